$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text values look like plain numbers need an explicit Text
# number format first, otherwise Excel auto-converts them to floating point
# numbers (losing exact text formatting / precision).
$ws.Range("D2").Value = "56.213.42"
$ws.Range("E2").Value = "  -5.03%  "
$ws.Range("D3").Value = "2.353.29"
$ws.Range("E3").Value = "  -6.65%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "509.68"
$ws.Range("E5").Value = "  -4.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.92"
$ws.Range("E6").Value = "  -5.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -2.77%  "
$ws.Range("D9").Value = "2.367.73"
$ws.Range("E9").Value = "  -6.03%  "
$ws.Range("E10").Value = "  -4.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.79"
$ws.Range("E12").Value = "  -8.55%  "
$ws.Range("E13").Value = "  -5.82%  "
$ws.Range("D14").Value = "2.773.94"
$ws.Range("E14").Value = "  -6.51%  "
$ws.Range("D15").Value = "56.160.61"
$ws.Range("E15").Value = "  -5.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.27"
$ws.Range("E16").Value = "  -5.25%  "
$ws.Range("E17").Value = "  -4.84%  "
$ws.Range("D18").Value = "2.329.85"
$ws.Range("E18").Value = "  -7.45%  "
$ws.Range("E20").Value = "  -4.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "308.98"
$ws.Range("E21").Value = "  -3.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.14"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.95"
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  -6.62%  "
$ws.Range("D27").Value = "2.467.31"
$ws.Range("E27").Value = "  -6.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.151"
$ws.Range("E28").Value = "  -5.69%  "
$ws.Range("E29").Value = "  -4.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.74"
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("E31").Value = "  -4.37%  "
$ws.Range("E32").Value = "  -7.35%  "
$ws.Range("E33").Value = "  -3.12%  "
$ws.Range("E34").Value = "  -8.20%  "
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.61"
$ws.Range("E37").Value = "  -3.11%  "
$ws.Range("E38").Value = "  -6.23%  "
$ws.Range("E39").Value = "  -7.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.52"
$ws.Range("E40").Value = "  -3.22%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("E42").Value = "  -7.43%  "
$ws.Range("E43").Value = "  -4.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.73"
$ws.Range("E44").Value = "  -7.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "122.94"
$ws.Range("E45").Value = "  -7.03%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "251.44"
$ws.Range("E46").Value = "  -9.75%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.564"
$ws.Range("E47").Value = "  -5.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0899"
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("E49").Value = "  -5.56%  "
$ws.Range("E50").Value = "  -6.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.49"
$ws.Range("E51").Value = "  -7.22%  "
